# CA2 demographic variables update
#
# Inserts two new numbered-list paragraphs ("Demographic Variables:" and
# "Geographic Metadata:") right after the existing "Economic Variables:"
# bullet, plus one extra blank paragraph, before the document's trailing
# blank paragraph / sectPr.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Locate the "Economic Variables:" list paragraph that the new content
# must follow, by searching rather than hard-coding an index.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Economic Variables:*") {
        $targetPara = $p
    }
}

# Paragraph 1: Demographic Variables
$demoXml = '<w:p ' + $wNs + '>' +
    '<w:pPr>' +
        '<w:numPr>' +
            '<w:ilvl w:val="0"/>' +
            '<w:numId w:val="3"/>' +
        '</w:numPr>' +
    '</w:pPr>' +
    '<w:r><w:t>Demographic Variables:</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Population growth rate and urbanization rate. These variables reflect societal trends, such as migration to urban areas, which can create increased housing demand in specific regions. A growing population in metropolitan areas often leads to higher property values due to increased competition for available housing.</w:t></w:r>' +
'</w:p>'

$null = $targetPara.Range.InsertParagraphAfter()
$demoPara = $targetPara.Next()
$null = $demoPara.Range.InsertXML($demoXml)

# Paragraph 2: Geographic Metadata
$geoXml = '<w:p ' + $wNs + '>' +
    '<w:pPr>' +
        '<w:numPr>' +
            '<w:ilvl w:val="0"/>' +
            '<w:numId w:val="3"/>' +
        '</w:numPr>' +
    '</w:pPr>' +
    '<w:r><w:t>Geographic Metadata:</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Country, city, region, and climate zone. This information allows for location-based analyses, enabling the identification of regional disparities and localized trends in housing prices. For instance, areas with </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>favorable</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> climates may experience higher housing prices due to their attractiveness to both residents and investors.</w:t></w:r>' +
'</w:p>'

$null = $demoPara.Range.InsertParagraphAfter()
$geoPara = $demoPara.Next()
$null = $geoPara.Range.InsertXML($geoXml)

# Extra blank paragraph inserted before the document's existing trailing
# blank paragraph. Built from raw XML (rather than relying on
# InsertParagraphAfter's formatting inheritance) so it comes out as a
# plain, empty <w:p/> with no inherited numbering/paragraph properties.
$null = $geoPara.Range.InsertParagraphAfter()
$blankPara = $geoPara.Next()
$null = $blankPara.Range.InsertXML('<w:p ' + $wNs + '/>')

Write-Output "CA2 demographic variables update applied."
